# Fruta / hortaliza, semanal
# Insert a new weekly record as row 183, shifting the existing rows
# (183..226) down by one (to 184..227).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(183).Insert()

$ws.Cells.Item(183, 1).Value  = 7
$ws.Cells.Item(183, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(183, 3).Value  = "Ñuble"
$ws.Cells.Item(183, 4).Value  = 44785
$ws.Cells.Item(183, 5).Value  = 16
$ws.Cells.Item(183, 6).Value  = 100112017
$ws.Cells.Item(183, 7).Value  = "Apio"
$ws.Cells.Item(183, 8).Value  = "Americana (o)"
$ws.Cells.Item(183, 9).Value  = "Primera"
$ws.Cells.Item(183, 10).Value = 120
$ws.Cells.Item(183, 11).Value = 9000
$ws.Cells.Item(183, 12).Value = 10000
$ws.Cells.Item(183, 13).Value = 9500
$ws.Cells.Item(183, 14).Value = "`$/docena de matas"
$ws.Cells.Item(183, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(183, 16).Value = 1583
$ws.Cells.Item(183, 17).Value = 6
$ws.Cells.Item(183, 18).Value = "Hortaliza"
